# Pdf letters and logger working
#
# Adds a new "Polling" column (H) to the worksheet with a header and three
# numeric values (matching rows 2-4 of the existing data), and makes sure
# the active selection / used range follow Excel's natural behaviour after
# such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in H1 (goes into the shared-string table alongside the rest).
$ws.Range("H1").Value = "Polling"

# New numeric data for the three existing data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 3

# Rename the built-in hyperlink cell style to its English name (no-op on
# hosts where the style collection doesn't support renaming).
try {
    $wb.Styles.Item(1).Name = "Hyperlink"
} catch {
}

# Leave the active cell on the last cell touched, like Excel would after
# typing the new values in order.
$null = $ws.Range("H4").Select()
